$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D:E values stay as text (avoid numeric auto-conversion of price strings)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '41.967.89'
$ws.Range("E2").Value = '  -4.49%  '
$ws.Range("D3").Value = '2.224.50'
$ws.Range("E3").Value = '  -5.44%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '244.23'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  -6.38%  '
$ws.Range("D7").Value = '68.52'
$ws.Range("E7").Value = '  -7.41%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  -8.61%  '
$ws.Range("E10").Value = '  -5.41%  '
$ws.Range("D11").Value = '58.21'
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("D12").Value = '35.44'
$ws.Range("E12").Value = '  +6.71%  '
$ws.Range("E13").Value = '  -3.51%  '
$ws.Range("E14").Value = '  -7.71%  '
$ws.Range("D15").Value = '2.554.61'
$ws.Range("E15").Value = '  -5.49%  '
$ws.Range("D16").Value = '14.77'
$ws.Range("E16").Value = '  -8.89%  '
$ws.Range("D17").Value = '0.842'
$ws.Range("E17").Value = '  -6.95%  '
$ws.Range("D18").Value = '2.225.90'
$ws.Range("E18").Value = '  -5.33%  '
$ws.Range("D19").Value = '41.906.79'
$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("E20").Value = '  -7.49%  '
$ws.Range("D21").Value = '72.48'
$ws.Range("E21").Value = '  -7.24%  '
$ws.Range("D22").Value = '6.12'
$ws.Range("E22").Value = '  -8.22%  '
$ws.Range("D23").Value = '234.59'
$ws.Range("E23").Value = '  -7.25%  '
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  +10.58%  '
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").Value = '3.62'
$ws.Range("E26").Value = '  -4.93%  '
$ws.Range("E27").Value = '  -2.44%  '
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").Value = '9.88'
$ws.Range("E29").Value = '  -5.64%  '
$ws.Range("D30").Value = '170.36'
$ws.Range("E30").Value = '  -3.50%  '
$ws.Range("E31").Value = '  -8.67%  '
$ws.Range("E32").Value = '  -6.22%  '
$ws.Range("E33").Value = '  -7.43%  '
$ws.Range("E34").Value = '  -5.08%  '
$ws.Range("D35").Value = '5.17'
$ws.Range("E35").Value = '  -4.18%  '
$ws.Range("D36").Value = '4.66'
$ws.Range("E36").Value = '  -8.33%  '
$ws.Range("D37").Value = '3.88'
$ws.Range("E37").Value = '  +1.54%  '
$ws.Range("D38").Value = '22.61'
$ws.Range("E38").Value = '  +17.77%  '
$ws.Range("D39").Value = '2.28'
$ws.Range("E39").Value = '  -4.83%  '
$ws.Range("D40").Value = '0.0276'
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("D41").Value = '5.83'
$ws.Range("E41").Value = '  -9.38%  '
$ws.Range("D42").Value = '66.36'
$ws.Range("E42").Value = '  +2.49%  '
$ws.Range("D43").Value = '4.97'
$ws.Range("E43").Value = '  -9.17%  '
$ws.Range("D44").Value = '8.93'
$ws.Range("E44").Value = '  -2.55%  '
$ws.Range("D45").Value = '0.101'
$ws.Range("E45").Value = '  -4.68%  '
$ws.Range("E46").Value = '  -5.41%  '
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").Value = '4.49'
$ws.Range("E48").Value = '  +4.98%  '
$ws.Range("B49").Value = 'Celestia'
$ws.Range("C49").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D49").Value = '10.17'
$ws.Range("E49").Value = '  +7.84%  '
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").Value = '1.18'
$ws.Range("E50").Value = '  -3.96%  '
$ws.Range("E51").Value = '  -4.62%  '

# Restore default cell style (remove temporary text-format style marker)
$ws.Range("D2:E51").Style = "Normal"

